$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 15243525
$ws.Range("I70").Value = 33534194
$ws.Range("J70").Value = 1300
$ws.Range("K70").Value = 100602582
$ws.Range("L70").Value = 3900
$ws.Range("M70").Value = -100602312
$ws.Range("N70").Value = -4440
$ws.Range("H73").Value = 15243525
$ws.Range("I73").Value = 33534194
$ws.Range("J73").Value = 1300
$ws.Range("K73").Value = 100602582
$ws.Range("L73").Value = 3900
$ws.Range("M73").Value = -100601646
$ws.Range("N73").Value = -5772
$ws.Range("H74").Value = 3889.9
$ws.Range("I74").Value = 3799.5
$ws.Range("J74").Value = 3912.5
$ws.Range("K74").Value = 3799.5
$ws.Range("L74").Value = 3912.5
$ws.Range("M74").Value = -2863.5
$ws.Range("N74").Value = -5784.5
$ws.Range("H77").Value = 3889.9
$ws.Range("I77").Value = 3799.5
$ws.Range("J77").Value = 3912.5
$ws.Range("K77").Value = 18997.5
$ws.Range("L77").Value = 19562.5
$ws.Range("M77").Value = -14317.5
$ws.Range("N77").Value = -28922.5
$ws.Range("H112").Value = 9572.866
$ws.Range("J112").Value = 12735.728
$ws.Range("L112").Value = 38207.18399999999
$ws.Range("N112").Value = -40423.18399999999
$ws.Range("H129").Value = 1032.0625
$ws.Range("J129").Value = 1099.4615
$ws.Range("L129").Value = 3298.3845
$ws.Range("N129").Value = -13298.3845
$ws.Range("H137").Value = 1933.3549
$ws.Range("I137").Value = 1644.2693
$ws.Range("J137").Value = 3436.6
$ws.Range("K137").Value = 4932.8079
$ws.Range("L137").Value = 10309.8
$ws.Range("M137").Value = -2382.8079
$ws.Range("N137").Value = -15409.8
$ws.Range("H138").Value = 4153.63
$ws.Range("J138").Value = 5091.415
$ws.Range("L138").Value = 15274.245
$ws.Range("N138").Value = -25554.245

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 57703.09
$ws.Range("I32").Value = 63455
$ws.Range("K32").Value = 63455
$ws.Range("M32").Value = -63168
$ws.Range("H63").Value = 3248.1667
$ws.Range("I63").Value = 2996.6667
$ws.Range("J63").Value = 3499.6667
$ws.Range("K63").Value = 2996.6667
$ws.Range("L63").Value = 3499.6667
$ws.Range("M63").Value = -2310.6667
$ws.Range("N63").Value = -4871.6667
$ws.Range("H66").Value = 3248.1667
$ws.Range("I66").Value = 2996.6667
$ws.Range("J66").Value = 3499.6667
$ws.Range("K66").Value = 14983.3335
$ws.Range("L66").Value = 17498.3335
$ws.Range("M66").Value = -11551.3335
$ws.Range("N66").Value = -24362.3335
$ws.Range("H74").Value = 1705.3182
$ws.Range("I74").Value = 1759.8823
$ws.Range("J74").Value = 1519.8
$ws.Range("K74").Value = 1759.8823
$ws.Range("L74").Value = 1519.8
$ws.Range("M74").Value = -885.8823
$ws.Range("N74").Value = -3267.8
$ws.Range("H77").Value = 1705.3182
$ws.Range("I77").Value = 1759.8823
$ws.Range("J77").Value = 1519.8
$ws.Range("K77").Value = 8799.4115
$ws.Range("L77").Value = 7599
$ws.Range("M77").Value = -4431.4115
$ws.Range("N77").Value = -16335
$ws.Range("H102").Value = 101641
$ws.Range("I102").Value = 1658.5714
$ws.Range("K102").Value = 1658.5714
$ws.Range("M102").Value = -36.57140000000004

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3756.3
$ws.Range("I105").Value = 3695.375
$ws.Range("K105").Value = 3695.375
$ws.Range("M105").Value = -1948.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4767.5
$ws.Range("J31").Value = 6160.1665
$ws.Range("L31").Value = 6160.1665
$ws.Range("N31").Value = -6750.1665
$ws.Range("H34").Value = 4767.5
$ws.Range("J34").Value = 6160.1665
$ws.Range("L34").Value = 6160.1665
$ws.Range("N34").Value = -6564.1665
$ws.Range("H99").Value = 3846.6667
$ws.Range("I99").Value = 3809.0908
$ws.Range("J99").Value = 3950
$ws.Range("K99").Value = 3809.0908
$ws.Range("L99").Value = 3950
$ws.Range("M99").Value = -2311.0908
$ws.Range("N99").Value = -6946
$ws.Range("H122").Value = 2024.4445
$ws.Range("I122").Value = 2132
$ws.Range("K122").Value = 6396
$ws.Range("M122").Value = -3946
$ws.Range("H126").Value = 3846.6667
$ws.Range("I126").Value = 3809.0908
$ws.Range("J126").Value = 3950
$ws.Range("K126").Value = 11427.2724
$ws.Range("L126").Value = 11850
$ws.Range("M126").Value = -8957.2724
$ws.Range("N126").Value = -16790

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3101.2554
$ws.Range("I113").Value = 640
$ws.Range("J113").Value = 3684.1843
$ws.Range("K113").Value = 1920
$ws.Range("L113").Value = 11052.5529
$ws.Range("M113").Value = 250
$ws.Range("N113").Value = -15392.5529
$ws.Range("H115").Value = 2307.261
$ws.Range("I115").Value = 277.33334
$ws.Range("J115").Value = 3023.7058
$ws.Range("K115").Value = 832.0000200000001
$ws.Range("L115").Value = 9071.117400000001
$ws.Range("M115").Value = 342.9999799999999
$ws.Range("N115").Value = -11421.1174
$ws.Range("H131").Value = 20003470
$ws.Range("I131").Value = 17208.334
$ws.Range("J131").Value = 22728868
$ws.Range("K131").Value = 51625.00199999999
$ws.Range("L131").Value = 68186604
$ws.Range("M131").Value = -46585.00199999999
$ws.Range("N131").Value = -68196684

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1412.5238
$ws.Range("I132").Value = 1220.5
$ws.Range("J132").Value = 2564.6667
$ws.Range("K132").Value = 3661.5
$ws.Range("L132").Value = 7694.000100000001
$ws.Range("M132").Value = -1131.5
$ws.Range("N132").Value = -12754.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 840.4761999999999
$ws.Range("I16").Value = 1159.4286
$ws.Range("J16").Value = 202.57143
$ws.Range("K16").Value = 1159.4286
$ws.Range("L16").Value = 202.57143
$ws.Range("M16").Value = -989.4286
$ws.Range("N16").Value = -542.57143
$ws.Range("H43").Value = 33975
$ws.Range("J43").Value = 17950
$ws.Range("L43").Value = 17950
$ws.Range("N43").Value = -18336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 91796.73
$ws.Range("I81").Value = 111973.336
$ws.Range("K81").Value = 223946.672
$ws.Range("M81").Value = -222885.672
$ws.Range("H84").Value = 91796.73
$ws.Range("I84").Value = 111973.336
$ws.Range("K84").Value = 1119733.36
$ws.Range("M84").Value = -1114429.36
